$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The paragraph right after the title currently just holds the
#    "_GoBack" bookmark (<w:bookmarkStart/><w:bookmarkEnd/>). Replace it
#    with a bare empty paragraph; the bookmark gets re-created later,
#    at the very end of the document.
# ---------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item(2)
$bookmarkRange = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.End)
$bookmarkRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# ---------------------------------------------------------------------
# 2) Append new content after the last paragraph ("in the response look
#    for refresh token"):
#       Links:
#       <hyperlink #1>
#       <hyperlink #2>
#       <empty paragraph holding the _GoBack bookmark>
# ---------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# "Links:" paragraph
$insPoint = $d.Range($d.Content.End, $d.Content.End)
$insPoint.InsertXML("<w:p $wNs><w:r><w:t>Links:</w:t></w:r></w:p>")

# First hyperlink paragraph
$insPoint = $d.Range($d.Content.End, $d.Content.End)
$insPoint.InsertXML("<w:p $wNs/>")
$para = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Hyperlinks.Add($para.Range, "https://www.dropboxforum.com/t5/Dropbox-API-Support-Feedback/Get-refresh-token-from-access-token/td-p/596739", "", "", "https://www.dropboxforum.com/t5/Dropbox-API-Support-Feedback/Get-refresh-token-from-access-token/td-p/596739") | Out-Null

# Second hyperlink paragraph
$insPoint = $d.Range($d.Content.End, $d.Content.End)
$insPoint.InsertXML("<w:p $wNs/>")
$para = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Hyperlinks.Add($para.Range, "https://www.dropbox.com/login?cont=https%3A%2F%2Fwww.dropbox.com%2Fdevelopers%2Fapps%3F_tk%3Dpilot_lp%26_ad%3Dtopbar4%26_camp%3Dmyapps", "", "", "https://www.dropbox.com/login?cont=https%3A%2F%2Fwww.dropbox.com%2Fdevelopers%2Fapps%3F_tk%3Dpilot_lp%26_ad%3Dtopbar4%26_camp%3Dmyapps") | Out-Null

# Trailing empty paragraph with the relocated "_GoBack" bookmark
$insPoint = $d.Range($d.Content.End, $d.Content.End)
$insPoint.InsertXML("<w:p $wNs><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>")
